$d = $word.ActiveDocument

# --- Locate the two target paragraphs -------------------------------------
# 1) The (currently empty) "Heading 2" paragraph that becomes the book
#    abbreviation title ("PHM") on the Philemon resource title page.
# 2) The very next paragraph (currently two empty runs) that gets a new
#    italic run listing every verse reference in the book.
$phmPara = $null
$versePara = $null

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($phmPara -eq $null) {
        # Paragraph.Range.Text carries a trailing paragraph-mark character
        # (CR, char 13) even when the paragraph is otherwise empty, so trim
        # it before comparing.
        $pText = $p.Range.Text.TrimEnd([char]13)
        if ($p.Style.NameLocal -eq "Heading 2" -and $pText -eq "") {
            $next = $d.Paragraphs.Item($i + 1)
            $nextText = $next.Range.Text.TrimEnd([char]13)
            if ($nextText -eq "") {
                $phmPara = $p
                $versePara = $next
            }
        }
    }
}

if ($phmPara -eq $null) {
    throw "Could not locate the empty 'Heading 2' title paragraph to fill in with the book abbreviation."
}

# --- 1) Give the heading paragraph its "PHM" run ---------------------------
# Insert at a collapsed (zero-length) range positioned at the start of the
# paragraph's existing (empty) run. A collapsed-range InsertXML keeps the
# paragraph's own formatting (pStyle Heading2) untouched and simply places
# the freshly-inserted run(s) ahead of whatever was already sitting at that
# position - which here is the single pre-existing empty run, giving us
# exactly: [new "PHM" run][original empty run].
$hStart = $phmPara.Range.Start
$hCollapsed = $d.Range($hStart, $hStart)

$phmXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t>PHM</w:t></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$hCollapsed.InsertXML($phmXml)

# --- 2) Add the italic verse-reference run to the following paragraph -----
# Build "Philemon 1:1, Philemon 1:2, ..., Philemon 1:25".
$verseText = ""
for ($v = 1; $v -le 25; $v++) {
    if ($v -gt 1) { $verseText += ", " }
    $verseText += "Philemon 1:$v"
}

# This paragraph already holds two empty runs sharing the same position
# (both zero-length), so a collapsed-range insert can't land "between" them.
# Instead we replace the whole (empty) paragraph range with the equivalent
# three runs: the original leading empty run, the new italic verse run,
# and the original trailing empty run - reproducing the diff's run order.
$vRange = $versePara.Range

$verseXml = @"
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr></w:r><w:r><w:rPr><w:i/><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr><w:t>$verseText</w:t></w:r><w:r><w:rPr><w:lang w:val="en_US" w:bidi="en_US"/></w:rPr></w:r></w:p><w:sectPr><w:pgSz w:w="12240" w:h="15840"/></w:sectPr></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
"@

$vRange.InsertXML($verseXml)
